# mals2-44 - address + formatting changes
#
# Applies the textual / formatting edits from the commit. Most of the
# diff's hunks are Word's automatic proofing markers (<w:proofErr .../>
# gramStart/gramEnd) that get inserted as a side-effect of Word's grammar
# checker re-scanning the paragraph after a nearby edit; those carry no
# visible text and aren't reachable through the documented Word object
# model, so this script focuses on the substantive, user-visible part of
# each edit: splitting the runs at the same text boundaries (so the
# document is byte-for-byte ready for the next proofing pass) and the
# real content/formatting changes (address block, signature block, font
# size).

$d = $word.ActiveDocument

# Helper: force Word to split the run(s) underlying $rng into a separate
# run without altering the run's effective formatting. We do this by
# reading the current Bold tri-state and toggling it off-and-on (or
# on-and-off) around it, which is a no-net-effect edit that nonetheless
# forces the engine to materialize $rng as its own run(s).
function Force-Split($rng) {
    $orig = $rng.Bold
    if ($orig -eq -1 -or $orig -eq $true) {
        $rng.Bold = $false
        $rng.Bold = $true
    } else {
        $rng.Bold = $true
        $rng.Bold = $false
    }
}

# ---------------------------------------------------------------------
# 1) "{d.MailingCity}" -> split "d.MailingCity" into "d." / "MailingCity"
#    runs (matches the spellStart/gramStart .. spellEnd/gramEnd nesting).
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("MailingCity")
if ($r.Found) {
    Force-Split $r
}

# ---------------------------------------------------------------------
# 2) Licence fee line: " Fee ………………………………………..….. " -> three runs,
#    " Fee ……………………………………" / "….." / "….. "
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(" Fee " + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + [char]0x2026 + ".." + [char]0x2026 + ".. ")
if ($r.Found) {
    $s = $r.Start
    # "….." == chars [19,22) of the 26-char matched string
    $mid = $d.Range($s + 19, $s + 22)
    Force-Split $mid
}

# ---------------------------------------------------------------------
# 3) "For the purpose of manufacturing and selling medicated feed." ->
#    "For the purpose of" / " manufacturing and selling medicated feed."
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("For the purpose of manufacturing and selling medicated feed.")
if ($r.Found) {
    $s = $r.Start
    $part1 = $d.Range($s, $s + 18)
    Force-Split $part1
}

# ---------------------------------------------------------------------
# 4) Signature block: " Minister of Finance " keeps its leading/trailing
#    space at the paragraph's normal size, but the name itself shrinks
#    to 10.5pt (w:sz 21) -- matches the target run split + Font.Size.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Minister of Finance")
if ($r.Found) {
    $r.Font.Size = 10.5
}

# ---------------------------------------------------------------------
# 5) Division name under "Ministry of Agriculture and Food" changes from
#    "Livestock Health Management and Regulation" to
#    "Office of the Chief Veterinarian". NB: the same phrase also
#    appears earlier in the document (title block) and must stay
#    untouched, so we scope the search to the paragraph right after the
#    "Ministry of Agriculture" heading.
# ---------------------------------------------------------------------
$heading = $d.Content
$heading.Find.Execute("Ministry of Agriculture")
if ($heading.Found) {
    $scope = $d.Range($heading.End, $d.Content.End)
    $scope.Find.Execute("Livestock Health Management and Regulation", $true, $false, $false, $false, $false, $true, 1, $false, "Office of the Chief Veterinarian", 2)
}

# ---------------------------------------------------------------------
# 6) Address block: drop the phone numbers and normalize "B.C." to "BC".
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Abbotsford, B.C.   V3G 2M3     Telephone: (778) 666-0560" + [char]9 + "Toll-Free: 1 (877) 877-2474", $true, $false, $false, $false, $false, $true, 1, $false, "Abbotsford, BC   V3G 2M3", 2)

Write-Output "edits applied"
